$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
Write-Host "Sheet1 name" $ws.Name
